$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: update the header/count values in B1:E1 ---
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# --- Row 2: drop the old B2/C2 pair, shift remaining values left ---
$ws.Range("B2").Value = 5.0175368920160865
$ws.Range("C2").ClearContents()
# D2 (7.1594031091868082) is unchanged
$ws.Range("E2").Value = 5.9453661633681243

# --- Row 3: drop the old B3 value, shift remaining values left ---
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 6.5446773553235111
$ws.Range("D3").Value = 6.1832668890764779
$ws.Range("E3").Value = 5.709867510217463
# F3 (5.7098675102174674) is unchanged

# --- Update the selection stored in the sheet view ---
$ws.Range("B1:E3").Select()
